$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Cells.Item(5, 12).Value = 0.5306122448979592
$ws.Cells.Item(5, 13).Value = 0.7323943661971831
$ws.Cells.Item(5, 14).Value = 0.6153846153846154
$ws.Cells.Item(5, 15).Value = 0.01375831231368951
$ws.Cells.Item(5, 16).Value = 0.08450704225352113
$ws.Cells.Item(5, 17).Value = 0.04038461538461535
$ws.Cells.Item(5, 18).Value = 0.02661934338952971
$ws.Cells.Item(5, 19).Value = 0.1304347826086956
$ws.Cells.Item(5, 20).Value = 0.07023411371237452

# Row 6
$ws.Cells.Item(6, 12).Value = 0.5306122448979592
$ws.Cells.Item(6, 13).Value = 0.7323943661971831
$ws.Cells.Item(6, 14).Value = 0.6153846153846154
$ws.Cells.Item(6, 15).Value = 0.01375831231368951
$ws.Cells.Item(6, 16).Value = 0.08450704225352113
$ws.Cells.Item(6, 17).Value = 0.04038461538461535
$ws.Cells.Item(6, 18).Value = 0.02661934338952971
$ws.Cells.Item(6, 19).Value = 0.1304347826086956
$ws.Cells.Item(6, 20).Value = 0.07023411371237452

# Row 10
$ws.Cells.Item(10, 12).Value = 0.5306122448979592
$ws.Cells.Item(10, 13).Value = 0.7323943661971831
$ws.Cells.Item(10, 14).Value = 0.6153846153846154
$ws.Cells.Item(10, 15).Value = 0.01375831231368951
$ws.Cells.Item(10, 16).Value = 0.08450704225352113
$ws.Cells.Item(10, 17).Value = 0.04038461538461535
$ws.Cells.Item(10, 18).Value = 0.02661934338952971
$ws.Cells.Item(10, 19).Value = 0.1304347826086956
$ws.Cells.Item(10, 20).Value = 0.07023411371237452

# Row 11
$ws.Cells.Item(11, 12).Value = 0.5333333333333333
$ws.Cells.Item(11, 13).Value = 0.7887323943661971
$ws.Cells.Item(11, 14).Value = 0.6363636363636364
$ws.Cells.Item(11, 15).Value = 0.01647940074906362
$ws.Cells.Item(11, 16).Value = 0.1408450704225351
$ws.Cells.Item(11, 17).Value = 0.06136363636363629
$ws.Cells.Item(11, 18).Value = 0.03188405797101439
$ws.Cells.Item(11, 19).Value = 0.217391304347826
$ws.Cells.Item(11, 20).Value = 0.1067193675889327

# Row 14
$ws.Cells.Item(14, 12).Value = 0.1972222222222222
$ws.Cells.Item(14, 13).Value = 1
$ws.Cells.Item(14, 14).Value = 0.3294663573085847
$ws.Cells.Item(14, 15).Value = 0.001691495965239015
$ws.Cells.Item(14, 16).Value = 0.01408450704225361
$ws.Cells.Item(14, 17).Value = 0.00312603096825842
$ws.Cells.Item(14, 18).Value = 0.00865079365079382
$ws.Cells.Item(14, 19).Value = 0.01428571428571438
$ws.Cells.Item(14, 20).Value = 0.009579052038449016

# Row 15
$ws.Cells.Item(15, 12).Value = 0.1977715877437326
$ws.Cells.Item(15, 13).Value = 1
$ws.Cells.Item(15, 14).Value = 0.3302325581395349
$ws.Cells.Item(15, 15).Value = 0.002240861486749401
$ws.Cells.Item(15, 16).Value = 0.01408450704225361
$ws.Cells.Item(15, 17).Value = 0.003892231799208612
$ws.Cells.Item(15, 18).Value = 0.01146040588937551
$ws.Cells.Item(15, 19).Value = 0.01428571428571438
$ws.Cells.Item(15, 20).Value = 0.01192691029900353

# Row 16
$ws.Cells.Item(16, 12).Value = 0.1955922865013774
$ws.Cells.Item(16, 13).Value = 1
$ws.Cells.Item(16, 14).Value = 0.3271889400921659
$ws.Cells.Item(16, 15).Value = 0.0000615602443942153687572726994403637946
$ws.Cells.Item(16, 16).Value = 0.01408450704225361
$ws.Cells.Item(16, 17).Value = 0.0008486137518395886
$ws.Cells.Item(16, 18).Value = 0.0003148366784732729
$ws.Cells.Item(16, 19).Value = 0.01428571428571438
$ws.Cells.Item(16, 20).Value = 0.002600394996708454

# Row 19
$ws.Cells.Item(19, 12).Value = 0.1972222222222222
$ws.Cells.Item(19, 13).Value = 1
$ws.Cells.Item(19, 14).Value = 0.3294663573085847
$ws.Cells.Item(19, 15).Value = 0.001691495965239015
$ws.Cells.Item(19, 16).Value = 0.01408450704225361
$ws.Cells.Item(19, 17).Value = 0.00312603096825842
$ws.Cells.Item(19, 18).Value = 0.00865079365079382
$ws.Cells.Item(19, 19).Value = 0.01428571428571438
$ws.Cells.Item(19, 20).Value = 0.009579052038449016

# Row 20
$ws.Cells.Item(20, 12).Value = 0.1972222222222222
$ws.Cells.Item(20, 13).Value = 1
$ws.Cells.Item(20, 14).Value = 0.3294663573085847
$ws.Cells.Item(20, 15).Value = 0.001691495965239015
$ws.Cells.Item(20, 16).Value = 0.01408450704225361
$ws.Cells.Item(20, 17).Value = 0.00312603096825842
$ws.Cells.Item(20, 18).Value = 0.00865079365079382
$ws.Cells.Item(20, 19).Value = 0.01428571428571438
$ws.Cells.Item(20, 20).Value = 0.009579052038449016

# Row 21
$ws.Cells.Item(21, 12).Value = 0.1955922865013774
$ws.Cells.Item(21, 13).Value = 1
$ws.Cells.Item(21, 14).Value = 0.3271889400921659
$ws.Cells.Item(21, 15).Value = 0.0000615602443942153687572726994403637946
$ws.Cells.Item(21, 16).Value = 0.01408450704225361
$ws.Cells.Item(21, 17).Value = 0.0008486137518395886
$ws.Cells.Item(21, 18).Value = 0.0003148366784732729
$ws.Cells.Item(21, 19).Value = 0.01428571428571438
$ws.Cells.Item(21, 20).Value = 0.002600394996708454

# Row 26
$ws.Cells.Item(26, 12).Value = 0.8571428571428571
$ws.Cells.Item(26, 13).Value = 0.1690140845070423
$ws.Cells.Item(26, 14).Value = 0.2823529411764706
$ws.Cells.Item(26, 15).Value = 0
$ws.Cells.Item(26, 16).Value = 0.0000000000000000555111512312578270211816
$ws.Cells.Item(26, 17).Value = 0
$ws.Cells.Item(26, 18).Value = 0
$ws.Cells.Item(26, 19).Value = 0.0000000000000003284409781182755915842064
$ws.Cells.Item(26, 20).Value = 0

# Row 29
$ws.Cells.Item(29, 12).Value = 0.8571428571428571
$ws.Cells.Item(29, 13).Value = 0.1690140845070423
$ws.Cells.Item(29, 14).Value = 0.2823529411764706
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 16).Value = 0.0000000000000000555111512312578270211816
$ws.Cells.Item(29, 17).Value = 0
$ws.Cells.Item(29, 18).Value = 0
$ws.Cells.Item(29, 19).Value = 0.0000000000000003284409781182755915842064
$ws.Cells.Item(29, 20).Value = 0

# Row 30
$ws.Cells.Item(30, 12).Value = 0.8571428571428571
$ws.Cells.Item(30, 13).Value = 0.1690140845070423
$ws.Cells.Item(30, 14).Value = 0.2823529411764706
$ws.Cells.Item(30, 15).Value = 0
$ws.Cells.Item(30, 16).Value = 0.0000000000000000555111512312578270211816
$ws.Cells.Item(30, 17).Value = 0
$ws.Cells.Item(30, 18).Value = 0
$ws.Cells.Item(30, 19).Value = 0.0000000000000003284409781182755915842064
$ws.Cells.Item(30, 20).Value = 0

# Row 31
$ws.Cells.Item(31, 12).Value = 0.8571428571428571
$ws.Cells.Item(31, 13).Value = 0.1690140845070423
$ws.Cells.Item(31, 14).Value = 0.2823529411764706
$ws.Cells.Item(31, 15).Value = 0
$ws.Cells.Item(31, 16).Value = 0.0000000000000000555111512312578270211816
$ws.Cells.Item(31, 17).Value = 0
$ws.Cells.Item(31, 18).Value = 0
$ws.Cells.Item(31, 19).Value = 0.0000000000000003284409781182755915842064
$ws.Cells.Item(31, 20).Value = 0

# Row 36
$ws.Cells.Item(36, 12).Value = 0.7352941176470589
$ws.Cells.Item(36, 13).Value = 0.352112676056338
$ws.Cells.Item(36, 14).Value = 0.4761904761904762
$ws.Cells.Item(36, 15).Value = 0.02100840336134457
$ws.Cells.Item(36, 16).Value = 0.2816901408450704
$ws.Cells.Item(36, 17).Value = 0.347985347985348
$ws.Cells.Item(36, 18).Value = 0.0294117647058824
$ws.Cells.Item(36, 19).Value = 4
$ws.Cells.Item(36, 20).Value = 2.714285714285714
